$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add a new row (row 5) to the change-track table, mirroring the
# formatting/style of the row above it (row 4) by copying row 4's
# formats down before filling in the new values.
$ws.Range("A4:E4").Copy()
$ws.Range("A5:E5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A5").Value = 41699

$ws.Range("B5").Value = "03"

$ws.Range("C5").Value = "JMR"

$ws.Range("D5").Value = "Update from GPT naming to PIT"

$ws.Range("E5").Value = "Done"

$ws.Range("E5").Select()
